$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-25: 45241 -> 45242
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45241) {
        $cell.Value = 45242
    }
}
